# Update automàtic: dades i banners [2026-02-11 18:20]
# Applies the per-station weather refresh captured in the commit diff:
# new extraction timestamps (col E) plus refreshed observation values
# (humidity, precipitation, pressure, wind gust, temperatures, etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = "'2026-02-11 18:18:28"
$ws.Range("H2").Value = "'79%"
$ws.Range("I2").Value = "'2.5 mm"

# Row 3
$ws.Range("E3").Value = "'2026-02-11 18:18:30"
$ws.Range("I3").Value = "'0.9 mm"
$ws.Range("O3").Value = "'0.2 °C"

# Row 4
$ws.Range("E4").Value = "'2026-02-11 18:18:33"
$ws.Range("H4").Value = "'56%"
$ws.Range("J4").Value = "'1002.7 hPa"

# Row 5
$ws.Range("E5").Value = "'2026-02-11 18:18:35"
$ws.Range("I5").Value = "'1.2 mm"
$ws.Range("O5").Value = "'0.5 °C"

# Row 6
$ws.Range("E6").Value = "'2026-02-11 18:18:38"
$ws.Range("J6").Value = "'1003.2 hPa"

# Row 7
$ws.Range("E7").Value = "'2026-02-11 18:18:40"
$ws.Range("J7").Value = "'1003.7 hPa"

# Row 8
$ws.Range("E8").Value = "'2026-02-11 18:18:43"
$ws.Range("N8").Value = "'13.3 °C 17:59 TU"

# Row 9
$ws.Range("E9").Value = "'2026-02-11 18:18:45"

# Row 10
$ws.Range("E10").Value = "'2026-02-11 18:18:48"
$ws.Range("H10").Value = "'76%"

# Row 11
$ws.Range("E11").Value = "'2026-02-11 18:18:53"
$ws.Range("O11").Value = "'8.0 °C"

# Row 12
$ws.Range("E12").Value = "'2026-02-11 18:18:55"

# Row 13
$ws.Range("E13").Value = "'2026-02-11 18:18:57"
$ws.Range("I13").Value = "'0.3 mm"
$ws.Range("O13").Value = "'7.5 °C"

# Row 14
$ws.Range("E14").Value = "'2026-02-11 18:19:00"
$ws.Range("O14").Value = "'19.2 °C"

# Row 15
$ws.Range("E15").Value = "'2026-02-11 18:19:02"

# Row 16
$ws.Range("E16").Value = "'2026-02-11 18:19:04"
$ws.Range("H16").Value = "'63%"
$ws.Range("I16").Value = "'5.9 mm"
$ws.Range("O16").Value = "'-0.6 °C"

# Row 17
$ws.Range("E17").Value = "'2026-02-11 18:19:07"

# Row 18
$ws.Range("E18").Value = "'2026-02-11 18:19:09"
$ws.Range("H18").Value = "'71%"
$ws.Range("J18").Value = "'1003.2 hPa"
$ws.Range("O18").Value = "'14.0 °C"

# Row 19
$ws.Range("E19").Value = "'2026-02-11 18:19:12"

# Row 20
$ws.Range("E20").Value = "'2026-02-11 18:19:14"
$ws.Range("I20").Value = "'0.8 mm"
$ws.Range("O20").Value = "'-1.1 °C"

# Row 21
$ws.Range("E21").Value = "'2026-02-11 18:19:17"
$ws.Range("I21").Value = "'1.7 mm"
$ws.Range("K21").Value = "'9.1 MJ/m2"
$ws.Range("O21").Value = "'8.4 °C"

# Row 22
$ws.Range("E22").Value = "'2026-02-11 18:19:24"
$ws.Range("H22").Value = "'91%"

# Row 23
$ws.Range("E23").Value = "'2026-02-11 18:19:26"
$ws.Range("H23").Value = "'71%"
$ws.Range("I23").Value = "'3.5 mm"

# Row 24
$ws.Range("E24").Value = "'2026-02-11 18:19:29"
$ws.Range("H24").Value = "'73%"
$ws.Range("I24").Value = "'6.8 mm"
$ws.Range("J24").Value = "'1007.2 hPa"
$ws.Range("N24").Value = "'11.1 °C 17:56 TU"
$ws.Range("O24").Value = "'13.4 °C"

# Row 25
$ws.Range("E25").Value = "'2026-02-11 18:19:31"
$ws.Range("H25").Value = "'62%"
$ws.Range("O25").Value = "'1.8 °C"

# Row 26
$ws.Range("E26").Value = "'2026-02-11 18:19:34"
$ws.Range("J26").Value = "'1003.1 hPa"

# Row 27
$ws.Range("E27").Value = "'2026-02-11 18:19:36"
$ws.Range("H27").Value = "'82%"
$ws.Range("I27").Value = "'1.1 mm"

# Row 28
$ws.Range("E28").Value = "'2026-02-11 18:19:39"
$ws.Range("J28").Value = "'1003.5 hPa"

# Row 29
$ws.Range("E29").Value = "'2026-02-11 18:19:41"

# Row 30
$ws.Range("E30").Value = "'2026-02-11 18:19:44"
$ws.Range("H30").Value = "'86%"
$ws.Range("J30").Value = "'1003.4 hPa"
$ws.Range("O30").Value = "'12.0 °C"

# Row 31
$ws.Range("E31").Value = "'2026-02-11 18:19:46"
$ws.Range("J31").Value = "'1002.6 hPa"

# Row 32
$ws.Range("E32").Value = "'2026-02-11 18:19:49"
$ws.Range("H32").Value = "'74%"

# Row 33
$ws.Range("E33").Value = "'2026-02-11 18:19:51"

# Row 34
$ws.Range("E34").Value = "'2026-02-11 18:19:54"
$ws.Range("H34").Value = "'60%"

# Row 35
$ws.Range("E35").Value = "'2026-02-11 18:19:56"
$ws.Range("H35").Value = "'69%"
$ws.Range("J35").Value = "'1008.0 hPa"

# Row 36
$ws.Range("E36").Value = "'2026-02-11 18:19:58"
$ws.Range("J36").Value = "'1003.5 hPa"
$ws.Range("O36").Value = "'13.0 °C"

# Row 37
$ws.Range("E37").Value = "'2026-02-11 18:20:01"
$ws.Range("H37").Value = "'81%"
$ws.Range("J37").Value = "'1004.6 hPa"
$ws.Range("L37").Value = "'38.2 km/h - 270º 17:48 TU"
$ws.Range("O37").Value = "'9.4 °C"

# Row 38
$ws.Range("E38").Value = "'2026-02-11 18:20:03"
$ws.Range("O38").Value = "'15.8 °C"

# Row 39
$ws.Range("E39").Value = "'2026-02-11 18:20:06"
$ws.Range("H39").Value = "'55%"
$ws.Range("L39").Value = "'94.7 km/h - 292º 17:52 TU"
$ws.Range("O39").Value = "'1.1 °C"

# Row 40
$ws.Range("E40").Value = "'2026-02-11 18:20:08"
$ws.Range("I40").Value = "'2.3 mm"
$ws.Range("J40").Value = "'1007.2 hPa"
$ws.Range("O40").Value = "'7.6 °C"

# Row 41
$ws.Range("E41").Value = "'2026-02-11 18:20:11"
$ws.Range("J41").Value = "'1005.0 hPa"
$ws.Range("N41").Value = "'15.9 °C 17:59 TU"

# Row 42
$ws.Range("E42").Value = "'2026-02-11 18:20:13"

# Row 43
$ws.Range("E43").Value = "'2026-02-11 18:20:16"

# Row 44
$ws.Range("E44").Value = "'2026-02-11 18:20:18"
$ws.Range("I44").Value = "'5.2 mm"

# Row 45
$ws.Range("E45").Value = "'2026-02-11 18:20:21"
$ws.Range("I45").Value = "'1.3 mm"
$ws.Range("J45").Value = "'1006.1 hPa"

# Row 46
$ws.Range("E46").Value = "'2026-02-11 18:20:23"
$ws.Range("H46").Value = "'57%"
$ws.Range("J46").Value = "'1007.6 hPa"
$ws.Range("N46").Value = "'13.0 °C 17:58 TU"
$ws.Range("O46").Value = "'17.4 °C"
